$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "mom1" label in H2 to "moment"
$ws.Range("H2").Value = "moment"

# Add two new columns: "mom" (I) and "pol" (J), matching the header style
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Value = "mom"

$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").Value = "pol"

# Add the data-row values for the new columns, matching the data-row style
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").Value = 1

$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").Value = 1

# Update the selected cell to match the saved selection state
$ws.Range("K13").Select() | Out-Null
